$d = $word.ActiveDocument

# The document ends with a ListParagraph item "Добавить отверстия для
# привинчивания циферблата." immediately followed (inside the same
# paragraph) by the "_GoBack" bookmark. We need to add a brand new
# ListParagraph bullet after it reading "Переместить разъем для
# батарейного отсека, чтобы было удобно паять провода." and have the
# "_GoBack" bookmark end up trailing the new bullet's text instead.
#
# Directly re-adding a collapsed bookmark at the very end of the
# document via Bookmarks.Add is unreliable here, so instead we let the
# existing bookmark's Range naturally carry itself forward by inserting
# our new sentence right after it (within the same run/paragraph), and
# only afterwards split the paragraph in two at the sentence boundary.
# That keeps the bookmark glued to the tail text, which becomes the new
# paragraph.

$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Range.InsertAfter("Переместить разъем для батарейного отсека, чтобы было удобно паять провода.")

# Locate the boundary between the two sentences (right after the
# existing "циферблата." text, right before the text we just added) and
# split the paragraph there so the new sentence becomes its own
# ListParagraph-styled bullet, carrying the bookmark along with it.
$boundary = $d.Content
$boundary.Find.Execute("циферблата.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$boundary.Collapse(0)
$boundary.InsertParagraphAfter()
